$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D cells whose new values parse as plain numbers need to be forced
# back to Text (matching the source data, which stores prices/¹ as strings)
# without leaving a residual custom style behind, so each such write is
# bracketed by a NumberFormat flip to "@" and a Style reset to "Normal".

# Row 2
$ws.Range("D2").Value = '42.829.64'
$ws.Range("E2").Value = '  -0.66%  '

# Row 3
$ws.Range("D3").Value = '2.326.82'
$ws.Range("E3").Value = '  -0.08%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.07%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.02'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.85%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '94.50'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.23%  '

# Row 7
$ws.Range("E7").Value = '  -0.25%  '

# Row 8
$ws.Range("E8").Value = '  -0.07%  '

# Row 9
$ws.Range("E9").Value = '  -2.14%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '33.99'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.42%  '

# Row 11
$ws.Range("B11").Value = 'Chainlink'
$ws.Range("C11").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '18.80'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.75%  '

# Row 12
$ws.Range("B12").Value = 'Dogecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0781'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.92%  '

# Row 13
$ws.Range("E13").Value = '  +2.08%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.71'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.47%  '

# Row 15
$ws.Range("D15").Value = '2.684.65'
$ws.Range("E15").Value = '  -0.25%  '

# Row 16
$ws.Range("D16").Value = '2.321.26'
$ws.Range("E16").Value = '  +0.18%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.789'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.17%  '

# Row 18
$ws.Range("D18").Value = '42.772.68'
$ws.Range("E18").Value = '  -0.61%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.96'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.19%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.16'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.54%  '

# Row 21
$ws.Range("D21").Value = '0.0₃0886'
$ws.Range("E21").Value = '  -1.55%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.88'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.11%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.36'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.93%  '

# Row 24
$ws.Range("E24").Value = '  +1.35%  '

# Row 25
$ws.Range("E25").Value = '  +0.00%  '

# Row 26
$ws.Range("E26").Value = '  -1.56%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.52'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.61%  '

# Row 28
$ws.Range("E28").Value = '  -0.69%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.11'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.33%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '31.29'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.92%  '

# Row 31
$ws.Range("E31").Value = '  -0.03%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '138.43'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -16.50%  '

# Row 33
$ws.Range("E33").Value = '  -0.37%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.55'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.87%  '

# Row 35
$ws.Range("E35").Value = '  -0.38%  '

# Row 36
$ws.Range("E36").Value = '  -0.94%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.35'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.22%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.82'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.94%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.101'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.36%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '22.41'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +23.61%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.74'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.43%  '

# Row 42
$ws.Range("E42").Value = '  -1.23%  '

# Row 43
$ws.Range("D43").Value = '1.933.91'
$ws.Range("E43").Value = '  -3.00%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0279'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.47%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.19'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.33%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.08'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.13%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.71'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.56%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.88'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.29%  '

# Row 49
$ws.Range("D49").Value = '2.552.60'
$ws.Range("E49").Value = '  -0.22%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '52.72'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.76%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.21'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.22%  '
